$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1852.8572
$ws.Range("J58").Value = 6250
$ws.Range("L58").Value = 18750
$ws.Range("N58").Value = -19050
$ws.Range("H70").Value = 3274.9167
$ws.Range("I70").Value = 1740
$ws.Range("J70").Value = 4371.2856
$ws.Range("K70").Value = 5220
$ws.Range("L70").Value = 13113.8568
$ws.Range("M70").Value = -4950
$ws.Range("N70").Value = -13653.8568
$ws.Range("H73").Value = 3274.9167
$ws.Range("I73").Value = 1740
$ws.Range("J73").Value = 4371.2856
$ws.Range("K73").Value = 5220
$ws.Range("L73").Value = 13113.8568
$ws.Range("M73").Value = -4284
$ws.Range("N73").Value = -14985.8568
$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13894866
$ws.Range("I32").Value = 14291691
$ws.Range("K32").Value = 14291691
$ws.Range("M32").Value = -14291404
$ws.Range("H61").Value = 2420.1072
$ws.Range("I61").Value = 2286.2083
$ws.Range("J61").Value = 3223.5
$ws.Range("K61").Value = 2286.2083
$ws.Range("L61").Value = 3223.5
$ws.Range("M61").Value = -2074.2083
$ws.Range("N61").Value = -3647.5
$ws.Range("H74").Value = 5030.5
$ws.Range("I74").Value = 5169.636
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 5169.636
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -4295.636
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 5030.5
$ws.Range("I77").Value = 5169.636
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 25848.18
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -21480.18
$ws.Range("N77").Value = -26236
$ws.Range("H132").Value = 2516.111
$ws.Range("I132").Value = 2540.32
$ws.Range("J132").Value = 2213.5
$ws.Range("K132").Value = 7620.960000000001
$ws.Range("L132").Value = 6640.5
$ws.Range("M132").Value = -5090.960000000001
$ws.Range("N132").Value = -11700.5
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 2420.1072
$ws.Range("I136").Value = 2286.2083
$ws.Range("J136").Value = 3223.5
$ws.Range("K136").Value = 6858.624899999999
$ws.Range("L136").Value = 9670.5
$ws.Range("M136").Value = -4308.624899999999
$ws.Range("N136").Value = -14770.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1242.6
$ws.Range("I20").Value = 1268.9166
$ws.Range("J20").Value = 1203.125
$ws.Range("K20").Value = 1268.9166
$ws.Range("L20").Value = 1203.125
$ws.Range("M20").Value = -1021.9166
$ws.Range("N20").Value = -1697.125
$ws.Range("H86").Value = 2723.3408
$ws.Range("J86").Value = 3897.4
$ws.Range("L86").Value = 3897.4
$ws.Range("N86").Value = -6143.4
$ws.Range("H89").Value = 2723.3408
$ws.Range("J89").Value = 3897.4
$ws.Range("L89").Value = 19487
$ws.Range("N89").Value = -30719
$ws.Range("H94").Value = 2313.0527
$ws.Range("I94").Value = 2536.5293
$ws.Range("K94").Value = 2536.5293
$ws.Range("M94").Value = -2085.5293
$ws.Range("H99").Value = 37974.184
$ws.Range("I99").Value = 45912.89
$ws.Range("K99").Value = 45912.89
$ws.Range("M99").Value = -44414.89
$ws.Range("H107").Value = 10161.375
$ws.Range("I107").Value = 8198.615
$ws.Range("K107").Value = 8198.615
$ws.Range("M107").Value = -6278.615
$ws.Range("H134").Value = 1717.0652
$ws.Range("I134").Value = 1362.4872
$ws.Range("J134").Value = 3692.5715
$ws.Range("K134").Value = 4087.4616
$ws.Range("L134").Value = 11077.7145
$ws.Range("M134").Value = -1552.4616
$ws.Range("N134").Value = -16147.7145

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2357.074
$ws.Range("J31").Value = 2356.3
$ws.Range("L31").Value = 2356.3
$ws.Range("N31").Value = -2946.3
$ws.Range("H34").Value = 2357.074
$ws.Range("J34").Value = 2356.3
$ws.Range("L34").Value = 2356.3
$ws.Range("N34").Value = -2760.3
$ws.Range("H134").Value = 2997.5334
$ws.Range("I134").Value = 2443.087
$ws.Range("J134").Value = 4819.2856
$ws.Range("K134").Value = 7329.261
$ws.Range("L134").Value = 14457.8568
$ws.Range("M134").Value = -4794.261
$ws.Range("N134").Value = -19527.8568
$ws.Range("H141").Value = 32571.072
$ws.Range("J141").Value = 32922.69
$ws.Range("L141").Value = 32922.69
$ws.Range("N141").Value = -43282.69

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5066.143
$ws.Range("J63").Value = 5599.8
$ws.Range("L63").Value = 16799.4
$ws.Range("N63").Value = -18297.4
$ws.Range("H66").Value = 5066.143
$ws.Range("J66").Value = 5599.8
$ws.Range("L66").Value = 50398.2
$ws.Range("N66").Value = -57886.2
$ws.Range("H107").Value = 294.75
$ws.Range("I107").Value = 326.66666
$ws.Range("K107").Value = 979.9999799999999
$ws.Range("M107").Value = 940.0000200000001
$ws.Range("H122").Value = 917.44446
$ws.Range("I122").Value = 1051.6
$ws.Range("J122").Value = 749.75
$ws.Range("K122").Value = 9464.4
$ws.Range("L122").Value = 6747.75
$ws.Range("M122").Value = -7014.4
$ws.Range("N122").Value = -11647.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7083.4736
$ws.Range("I5").Value = 7083.4736
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 7083.4736
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("M5").Value = -6971.4736
$ws.Range("H9").Value = 1805
$ws.Range("I9").Value = 2600.5
$ws.Range("J9").Value = 1274.6666
$ws.Range("K9").Value = 2600.5
$ws.Range("L9").Value = 1274.6666
$ws.Range("M9").Value = -2430.5
$ws.Range("N9").Value = -1614.6666
$ws.Range("H52").Value = 22397.6
$ws.Range("J52").Value = 26499.5
$ws.Range("L52").Value = 26499.5
$ws.Range("N52").Value = -27017.5
$ws.Range("H102").Value = 6021.3335
$ws.Range("I102").Value = 3806.4443
$ws.Range("K102").Value = 3806.4443
$ws.Range("M102").Value = -2184.4443
$ws.Range("H107").Value = 1027.75
$ws.Range("H122").Value = 8283.143
$ws.Range("I122").Value = 8164
$ws.Range("K122").Value = 24492
$ws.Range("M122").Value = -22042

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12866.9
$ws.Range("I40").Value = 14814.25
$ws.Range("J40").Value = 5077.5
$ws.Range("K40").Value = 14814.25
$ws.Range("L40").Value = 5077.5
$ws.Range("M40").Value = -14678.25
$ws.Range("N40").Value = -5349.5
$ws.Range("H61").Value = 2819.5
$ws.Range("I61").Value = 2206.3333
$ws.Range("K61").Value = 2206.3333
$ws.Range("M61").Value = -2004.3333
$ws.Range("H68").Value = 3817.375
$ws.Range("I68").Value = 3768.6667
$ws.Range("K68").Value = 3768.6667
$ws.Range("M68").Value = -3019.6667
$ws.Range("H71").Value = 3817.375
$ws.Range("I71").Value = 3768.6667
$ws.Range("K71").Value = 18843.3335
$ws.Range("M71").Value = -15099.3335
$ws.Range("H82").Value = 1102
$ws.Range("I82").Value = 1332.2307
$ws.Range("J82").Value = 902.4666999999999
$ws.Range("K82").Value = 1332.2307
$ws.Range("L82").Value = 902.4666999999999
$ws.Range("M82").Value = -971.2307000000001
$ws.Range("N82").Value = -1624.4667
$ws.Range("H85").Value = 1102
$ws.Range("I85").Value = 1332.2307
$ws.Range("J85").Value = 902.4666999999999
$ws.Range("K85").Value = 1332.2307
$ws.Range("L85").Value = 902.4666999999999
$ws.Range("M85").Value = -84.23070000000007
$ws.Range("N85").Value = -3398.4667
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H113").Value = 2819.5
$ws.Range("I113").Value = 2206.3333
$ws.Range("K113").Value = 2206.3333
$ws.Range("M113").Value = -36.33329999999978
$ws.Range("H122").Value = 5622.5454
$ws.Range("I122").Value = 2536.375
$ws.Range("J122").Value = 7386.0713
$ws.Range("K122").Value = 7609.125
$ws.Range("L122").Value = 22158.2139
$ws.Range("M122").Value = -5159.125
$ws.Range("N122").Value = -27058.2139

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40630
$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -42184
$ws.Range("H81").Value = 9262956
$ws.Range("J81").Value = 18523318
$ws.Range("L81").Value = 37046636
$ws.Range("N81").Value = -37048758
$ws.Range("H84").Value = 9262956
$ws.Range("J84").Value = 18523318
$ws.Range("L84").Value = 185233180
$ws.Range("N84").Value = -185243788
